$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 54, shifting existing rows 54-127 down to 55-128.
$ws.Rows.Item(54).Insert()

# Populate the newly inserted row 54 with a new data record
# (same market/category info as the row that used to be there, but with a
# new date and a new "Volumen" value).
$row = 54
$ws.Cells.Item($row, 1).Value = 8
$ws.Cells.Item($row, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item($row, 3).Value = "Coquimbo"
$ws.Cells.Item($row, 4).Value = 44413
$ws.Cells.Item($row, 5).Value = 4
$ws.Cells.Item($row, 6).Value = 100112012
$ws.Cells.Item($row, 7).Value = "Espinaca"
$ws.Cells.Item($row, 8).Value = "Sin especificar"
$ws.Cells.Item($row, 9).Value = "Primera"
$ws.Cells.Item($row, 10).Value = 3140
$ws.Cells.Item($row, 11).Value = 400
$ws.Cells.Item($row, 12).Value = 500
$ws.Cells.Item($row, 13).Value = 450
$ws.Cells.Item($row, 14).Value = "`$/atado 300 a 500 gramos"
$ws.Cells.Item($row, 15).Value = "Provincia del Elquí"
$ws.Cells.Item($row, 16).Value = 900
$ws.Cells.Item($row, 17).Value = 0.5
$ws.Cells.Item($row, 18).Value = "Hortaliza"
